$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 674396.1
$ws.Range("J129").Value = 789112.4399999999
$ws.Range("L129").Value = 2367337.32
$ws.Range("N129").Value = -2377337.32
$ws.Range("H132").Value = 2327976.2
$ws.Range("I132").Value = 2308.625
$ws.Range("K132").Value = 6925.875
$ws.Range("M132").Value = -4395.875
$ws.Range("H137").Value = 1060.7273
$ws.Range("I137").Value = 940.8889
$ws.Range("J137").Value = 1600
$ws.Range("K137").Value = 2822.6667
$ws.Range("L137").Value = 4800
$ws.Range("M137").Value = -272.6667000000002
$ws.Range("N137").Value = -9900
$ws.Range("H138").Value = 4183.352
$ws.Range("I138").Value = 2264.3684
$ws.Range("J138").Value = 4711.768
$ws.Range("K138").Value = 6793.1052
$ws.Range("L138").Value = 14135.304
$ws.Range("M138").Value = -1653.1052
$ws.Range("N138").Value = -24415.304

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1853.5
$ws.Range("I61").Value = 1611.4615
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1611.4615
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1399.4615
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 1345.9
$ws.Range("I74").Value = 1385.3914
$ws.Range("J74").Value = 1216.1428
$ws.Range("K74").Value = 1385.3914
$ws.Range("L74").Value = 1216.1428
$ws.Range("M74").Value = -511.3914
$ws.Range("N74").Value = -2964.1428
$ws.Range("H77").Value = 1345.9
$ws.Range("I77").Value = 1385.3914
$ws.Range("J77").Value = 1216.1428
$ws.Range("K77").Value = 6926.957
$ws.Range("L77").Value = 6080.714
$ws.Range("M77").Value = -2558.957
$ws.Range("N77").Value = -14816.714
$ws.Range("H136").Value = 1853.5
$ws.Range("I136").Value = 1611.4615
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4834.3845
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2284.3845
$ws.Range("N136").Value = -20100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20820.396
$ws.Range("I134").Value = 1663.2046
$ws.Range("J134").Value = 114477.78
$ws.Range("K134").Value = 4989.6138
$ws.Range("L134").Value = 343433.34
$ws.Range("M134").Value = -2454.6138
$ws.Range("N134").Value = -348503.34

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4797.3887
$ws.Range("I31").Value = 5189.6
$ws.Range("J31").Value = 4307.125
$ws.Range("K31").Value = 5189.6
$ws.Range("L31").Value = 4307.125
$ws.Range("M31").Value = -4894.6
$ws.Range("N31").Value = -4897.125
$ws.Range("H34").Value = 4797.3887
$ws.Range("I34").Value = 5189.6
$ws.Range("J34").Value = 4307.125
$ws.Range("K34").Value = 5189.6
$ws.Range("L34").Value = 4307.125
$ws.Range("M34").Value = -4987.6
$ws.Range("N34").Value = -4711.125
$ws.Range("H58").Value = 916.44446
$ws.Range("I58").Value = 843.4838999999999
$ws.Range("J58").Value = 1368.8
$ws.Range("K58").Value = 843.4838999999999
$ws.Range("L58").Value = 1368.8
$ws.Range("M58").Value = -640.4838999999999
$ws.Range("N58").Value = -1774.8
$ws.Range("H107").Value = 293.3125
$ws.Range("I107").Value = 145.66667
$ws.Range("J107").Value = 381.9
$ws.Range("K107").Value = 145.66667
$ws.Range("L107").Value = 381.9
$ws.Range("M107").Value = 1774.33333
$ws.Range("N107").Value = -4221.9
$ws.Range("H132").Value = 1350.1143
$ws.Range("I132").Value = 1008.5926
$ws.Range("K132").Value = 3025.7778
$ws.Range("M132").Value = -495.7777999999998
$ws.Range("H133").Value = 20000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 13514671
$ws.Range("I134").Value = 1200.3871
$ws.Range("J134").Value = 83334270
$ws.Range("K134").Value = 3601.1613
$ws.Range("L134").Value = 250002810
$ws.Range("M134").Value = -1066.1613
$ws.Range("N134").Value = -250007880
$ws.Range("H135").Value = 56350
$ws.Range("J135").Value = 56350
$ws.Range("L135").Value = 56350
$ws.Range("N135").Value = -66490
$ws.Range("H136").Value = 916.44446
$ws.Range("I136").Value = 843.4838999999999
$ws.Range("J136").Value = 1368.8
$ws.Range("K136").Value = 2530.4517
$ws.Range("L136").Value = 4106.4
$ws.Range("M136").Value = 19.54830000000038
$ws.Range("N136").Value = -9206.4
$ws.Range("H137").Value = 79780
$ws.Range("J137").Value = 79780
$ws.Range("L137").Value = 79780
$ws.Range("N137").Value = -89980
$ws.Range("H138").Value = 77780
$ws.Range("J138").Value = 77780
$ws.Range("L138").Value = 77780
$ws.Range("N138").Value = -88060
$ws.Range("N133").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 13052
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 15065
$ws.Range("K110").Value = 15000
$ws.Range("L110").Value = 45195
$ws.Range("M110").Value = -10910
$ws.Range("N110").Value = -53375
$ws.Range("H131").Value = 20918010
$ws.Range("J131").Value = 64273.188
$ws.Range("L131").Value = 192819.564
$ws.Range("N131").Value = -202899.564

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3261.9092
$ws.Range("I132").Value = 3250.875
$ws.Range("J132").Value = 3291.3333
$ws.Range("K132").Value = 9752.625
$ws.Range("L132").Value = 9873.999899999999
$ws.Range("M132").Value = -7222.625
$ws.Range("N132").Value = -14933.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2933.3333
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2933.3333
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2933.3333
$ws.Range("N68").Value = -4431.3333
$ws.Range("H71").Value = 2933.3333
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2933.3333
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14666.6665
$ws.Range("N71").Value = -22154.6665
$ws.Range("H136").Value = 3447.138
$ws.Range("I136").Value = 1988.9048
$ws.Range("J136").Value = 7275
$ws.Range("K136").Value = 5966.7144
$ws.Range("L136").Value = 21825
$ws.Range("M136").Value = -3416.7144
$ws.Range("N136").Value = -26925
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5124.48
$ws.Range("I136").Value = 953.86664
$ws.Range("J136").Value = 11380.4
$ws.Range("K136").Value = 2861.59992
$ws.Range("L136").Value = 34141.2
$ws.Range("M136").Value = -311.5999199999997
$ws.Range("N136").Value = -39241.2
